$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts existing rows 13+ down by one),
# matching Excel's default "insert copies format from the row above" behaviour.
$ws.Rows("13:13").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Re-apply the exact formatting used by the row above (row 12) so the new
# row's cell styles line up with the rest of the R_ROLE_PERM block.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# New permission row: same roleId as every other row in this block, a new
# permId, and the Chinese label for the new "update environment info" API.
$ws.Range("A13").Value = "e501b47a-c08b-4c83-b12b-95ad82873e96"
$ws.Range("B13").Value = "62e52251-6e19-4ec5-a7d3-cfdf2968d4ca"
$ws.Range("C13").Value = "环境信息更新"

# The selected band shifts from rows 13:15 to 15:15 after the insert.
$ws.Range("A15:XFD15").Select() | Out-Null
